# TestExperiment4.xlsx - "CHISA Fighting" commit: update a couple of
# experiment-condition inputs and correct a value that had been entered as
# text back into a proper number, then leave the selection where the user
# last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Flowrate (B2) and Feed volume (B3) experiment-condition inputs were updated.
$ws.Range("B2").Value = 235
$ws.Range("B3").Value = 16

# E25 had been typed in as the text "0,9417" (a shared string); fix it up as
# the numeric value 0.9417 so it behaves like the other data-table cells.
$ws.Range("E25").Value = 0.9417

# Leave the selection on G12, matching where the user ended up.
$ws.Range("G12").Select()
